$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column L header ---
$ws.Range("L1").Value = "经销商索赔代码-经销商索赔单号"

# --- Update dealer claim code column (B) for all data rows ---
$ws.Range("B2").Value = "LB144870518080110"
$ws.Range("B3").Value = "LB144870518080110"
$ws.Range("B4").Value = "LB144870518080110"

# --- Update claim-order numbers in column E ---
$ws.Range("E2").Value = 1000054
$ws.Range("E3").Value = 1000055
$ws.Range("E4").Value = 3001232

# Copy the formatting used by the other green "text" cells (column D) onto
# the new column L cells so they pick up the same style (numFmt "@" + fill).
$ws.Range("D2").Copy() | Out-Null
$ws.Range("L2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("D3").Copy() | Out-Null
$ws.Range("L3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("D4").Copy() | Out-Null
$ws.Range("L4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# --- New column L value (row 2 only; rows 3-4 stay blank but styled) ---
$ws.Range("L2").Value = "99956-1000055,99956-3001232"

# --- Column width for the new column ---
$ws.Columns.Item(12).ColumnWidth = 25.33

# --- Final selection matches the author's last cursor position ---
$ws.Range("L8").Select() | Out-Null
